$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new I/J columns (rows 1-65): a third parametric curve driven by
# column G with a /16 divisor (half the /32 divisor used by E/F).
for ($r = 1; $r -le 65; $r++) {
    $ws.Cells.Item($r, 9).Formula  = "=150-30*COS((P$r/2)*G$r/16)"
    $ws.Cells.Item($r, 10).Formula = "=70-30*SIN((P$r/2)*G$r/16)"
}

# G49 was accidentally typed as 16 instead of -16; fix it.
$ws.Range("G49").Value = -16

# Restore the view/selection as last saved by the author.
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("G50").Select()
